# "added delay in Get_Flight_Info"
#
# The flight that used to be tracked as row 6 (American, on-time) is now
# delayed, so a new "Spirit" flight entry (with a long, overnight delay)
# takes its place in row 6, and the previous row 7 flight (United) is
# replaced by the American flight's original info, shifted down to row 7.
#
# Net effect on the "Airline" mini-table (rows 5-7):
#   Row 6: American -> Spirit / new schedule / new (much longer) duration / new cost
#   Row 7: United    -> American (prior row 6 schedule/duration) / new cost

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 becomes what row 6 used to be (American flight), with an updated ticket cost.
$ws.Range("A7").Value2 = "American"
$ws.Range("B7").Value2 = "5:02 AM – 10:05 AM"
$ws.Range("D7").Value2 = "7 hr 3 min"
$ws.Range("E7").Value2 = 507

# Row 6 becomes the new, delayed Spirit flight.
$ws.Range("A6").Value2 = "Spirit"
$ws.Range("B6").Value2 = "10:50 AM – 7:02 PM+1"
$ws.Range("D6").Value2 = "34 hr 12 min"
$ws.Range("E6").Value2 = 414
